$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1312.3846
$ws.Range("I19").Value = 1269.7142
$ws.Range("J19").Value = 1362.1666
$ws.Range("K19").Value = 1269.7142
$ws.Range("L19").Value = 1362.1666
$ws.Range("M19").Value = -1094.7142
$ws.Range("N19").Value = -1712.1666
$ws.Range("H33").Value = 467.13635
$ws.Range("I33").Value = 292.86667
$ws.Range("J33").Value = 840.5714
$ws.Range("K33").Value = 292.86667
$ws.Range("L33").Value = 840.5714
$ws.Range("M33").Value = -63.86667
$ws.Range("N33").Value = -1298.5714
$ws.Range("H113").Value = 9367.786
$ws.Range("I113").Value = 14466.125
$ws.Range("J113").Value = 2570
$ws.Range("K113").Value = 14466.125
$ws.Range("L113").Value = 2570
$ws.Range("M113").Value = -11212.125
$ws.Range("N113").Value = -9078
$ws.Range("H132").Value = 6158.387
$ws.Range("I132").Value = 6197
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 18591
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -16061
$ws.Range("N132").Value = -20060
$ws.Range("H135").Value = 428.44827
$ws.Range("I135").Value = 425.89285
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 3833.03565
$ws.Range("L135").Value = 4500
$ws.Range("M135").Value = -1298.03565
$ws.Range("N135").Value = -9570
$ws.Range("H137").Value = 27360.441
$ws.Range("I137").Value = 31434.31
$ws.Range("J137").Value = 3732
$ws.Range("K137").Value = 94302.93000000001
$ws.Range("L137").Value = 11196
$ws.Range("M137").Value = -91752.93000000001
$ws.Range("N137").Value = -16296
$ws.Range("H138").Value = 2536.81
$ws.Range("I138").Value = 856.9091
$ws.Range("J138").Value = 3856.7322
$ws.Range("K138").Value = 2570.7273
$ws.Range("L138").Value = 11570.1966
$ws.Range("M138").Value = 2569.2727
$ws.Range("N138").Value = -21850.1966

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16783474
$ws.Range("I32").Value = 14287698
$ws.Range("J32").Value = 83337496
$ws.Range("K32").Value = 14287698
$ws.Range("L32").Value = 83337496
$ws.Range("M32").Value = -14287411
$ws.Range("N32").Value = -83338070
$ws.Range("H61").Value = 2485.303
$ws.Range("I61").Value = 2305.3794
$ws.Range("J61").Value = 3789.75
$ws.Range("K61").Value = 2305.3794
$ws.Range("L61").Value = 3789.75
$ws.Range("M61").Value = -2093.3794
$ws.Range("N61").Value = -4213.75
$ws.Range("H74").Value = 2033.1698
$ws.Range("I74").Value = 1909.16
$ws.Range("J74").Value = 4100
$ws.Range("K74").Value = 1909.16
$ws.Range("L74").Value = 4100
$ws.Range("M74").Value = -1035.16
$ws.Range("N74").Value = -5848
$ws.Range("H77").Value = 2033.1698
$ws.Range("I77").Value = 1909.16
$ws.Range("J77").Value = 4100
$ws.Range("K77").Value = 9545.800000000001
$ws.Range("L77").Value = 20500
$ws.Range("M77").Value = -5177.800000000001
$ws.Range("N77").Value = -29236
$ws.Range("H109").Value = 65827.836
$ws.Range("J109").Value = 65827.836
$ws.Range("L109").Value = 65827.836
$ws.Range("N109").Value = -68601.836
$ws.Range("H112").Value = 43064
$ws.Range("J112").Value = 43064
$ws.Range("L112").Value = 43064
$ws.Range("N112").Value = -46018
$ws.Range("H132").Value = 139051.06
$ws.Range("I132").Value = 162980.81
$ws.Range("K132").Value = 488942.43
$ws.Range("M132").Value = -486412.43
$ws.Range("H136").Value = 2485.303
$ws.Range("I136").Value = 2305.3794
$ws.Range("J136").Value = 3789.75
$ws.Range("K136").Value = 6916.138199999999
$ws.Range("L136").Value = 11369.25
$ws.Range("M136").Value = -4366.138199999999
$ws.Range("N136").Value = -16469.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1719.0857
$ws.Range("I94").Value = 1394.96
$ws.Range("J94").Value = 2529.4
$ws.Range("K94").Value = 1394.96
$ws.Range("L94").Value = 2529.4
$ws.Range("M94").Value = -943.96
$ws.Range("N94").Value = -3431.4
$ws.Range("H105").Value = 2271.682
$ws.Range("I105").Value = 1798.6154
$ws.Range("J105").Value = 2955
$ws.Range("K105").Value = 1798.6154
$ws.Range("L105").Value = 2955
$ws.Range("M105").Value = -51.61539999999991
$ws.Range("N105").Value = -6449
$ws.Range("H110").Value = 35000.5
$ws.Range("J110").Value = 35000.5
$ws.Range("L110").Value = 35000.5
$ws.Range("N110").Value = -43180.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2622.8157
$ws.Range("I31").Value = 1922.7142
$ws.Range("K31").Value = 1922.7142
$ws.Range("M31").Value = -1627.7142
$ws.Range("H34").Value = 2622.8157
$ws.Range("I34").Value = 1922.7142
$ws.Range("K34").Value = 1922.7142
$ws.Range("M34").Value = -1720.7142
$ws.Range("H58").Value = 2105.7222
$ws.Range("I58").Value = 1838.64
$ws.Range("J58").Value = 2712.7273
$ws.Range("K58").Value = 1838.64
$ws.Range("L58").Value = 2712.7273
$ws.Range("M58").Value = -1635.64
$ws.Range("N58").Value = -3118.7273
$ws.Range("H134").Value = 1779.1714
$ws.Range("I134").Value = 1605.5161
$ws.Range("J134").Value = 3125
$ws.Range("K134").Value = 4816.5483
$ws.Range("L134").Value = 9375
$ws.Range("M134").Value = -2281.5483
$ws.Range("N134").Value = -14445
$ws.Range("H136").Value = 2105.7222
$ws.Range("I136").Value = 1838.64
$ws.Range("J136").Value = 2712.7273
$ws.Range("K136").Value = 5515.92
$ws.Range("L136").Value = 8138.1819
$ws.Range("M136").Value = -2965.92
$ws.Range("N136").Value = -13238.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1350.3636
$ws.Range("I132").Value = 1236.7142
$ws.Range("K132").Value = 11130.4278
$ws.Range("M132").Value = -8600.427799999999
$ws.Range("H141").Value = 7998.5713
$ws.Range("I141").Value = 6198
$ws.Range("J141").Value = 12500
$ws.Range("K141").Value = 18594
$ws.Range("L141").Value = 37500
$ws.Range("M141").Value = -13414
$ws.Range("N141").Value = -47860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 74997
$ws.Range("J111").Value = 74997
$ws.Range("L111").Value = 74997
$ws.Range("N111").Value = -81131
$ws.Range("H113").Value = 3039.4285
$ws.Range("I113").Value = 2450
$ws.Range("J113").Value = 4513
$ws.Range("K113").Value = 2450
$ws.Range("L113").Value = 4513
$ws.Range("M113").Value = -280
$ws.Range("N113").Value = -8853
$ws.Range("H126").Value = 3174.111
$ws.Range("I126").Value = 3660.4
$ws.Range("K126").Value = 10981.2
$ws.Range("M126").Value = -8511.200000000001
$ws.Range("H132").Value = 2008.0476
$ws.Range("I132").Value = 1697.7858
$ws.Range("K132").Value = 5093.357400000001
$ws.Range("M132").Value = -2563.357400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 55733.25
$ws.Range("J110").Value = 55733.25
$ws.Range("L110").Value = 55733.25
$ws.Range("N110").Value = -63913.25
$ws.Range("H132").Value = 911253.75
$ws.Range("I132").Value = 1251724
$ws.Range("K132").Value = 3755172
$ws.Range("M132").Value = -3752642

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26311.186
$ws.Range("I132").Value = 33169.273
$ws.Range("J132").Value = 3679.5
$ws.Range("K132").Value = 99507.819
$ws.Range("L132").Value = 11038.5
$ws.Range("M132").Value = -96977.819
$ws.Range("N132").Value = -16098.5
